# Regenerate the handback report: the 64fac99d-... file now failed its
# handback transform, so its row moves up to directly follow the
# "Handed back" rows (row 6) instead of sitting after 70058cc9 / 4dd38fef,
# and its status flips from "Ready for handoff" to "Handback transform
# failed". The two rows it displaces (70058cc9, 4dd38fef) simply shift
# down to rows 7/8, keeping their own data untouched.

$wb = $excel.ActiveWorkbook

function Set-LinkText($ws, [string]$addr, [string]$text) {
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $addr) {
            $h.TextToDisplay = $text
        }
    }
}

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A6").Value = "64fac99d-8603-4021-9534-a43a7f1cfd1d.md"
$ws1.Range("B6:C6").Value = "Handback transform failed"

$ws1.Range("A7").Value = "70058cc9-c335-4816-b6d1-eb41c2e484ea.md"
$ws1.Range("B7:C7").Value = "In Translation"

$ws1.Range("A8").Value = "4dd38fef-441a-4be0-b79a-ef6c9247ebcf.md"
$ws1.Range("B8:C8").Value = "Ready for handoff"

Set-LinkText $ws1 '$A$6' "64fac99d-8603-4021-9534-a43a7f1cfd1d.md"
Set-LinkText $ws1 '$A$7' "70058cc9-c335-4816-b6d1-eb41c2e484ea.md"
Set-LinkText $ws1 '$A$8' "4dd38fef-441a-4be0-b79a-ef6c9247ebcf.md"

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A6").Value = "64fac99d-8603-4021-9534-a43a7f1cfd1d.md"
$ws2.Range("B6").Value = "Handback transform failed"
$ws2.Range("C6").Value = "64fac99d-8603-4021-9534-a43a7f1cfd1d.efd74baa64930741992ee27e4b2834cc8ef1667d.zh-cn.xlf"
$ws2.Range("D6").Value = "2016-03-10 04:20:43"

$ws2.Range("A7").Value = "70058cc9-c335-4816-b6d1-eb41c2e484ea.md"
$ws2.Range("B7").Value = "In Translation"
$ws2.Range("C7").Value = "70058cc9-c335-4816-b6d1-eb41c2e484ea.d5899fcb1515f857962642b1ceab8a68295ec2fb.zh-cn.xlf"
$ws2.Range("D7").Value = "2016-03-10 04:13:01"

$ws2.Range("A8").Value = "4dd38fef-441a-4be0-b79a-ef6c9247ebcf.md"
$ws2.Range("B8").Value = "Ready for handoff"
$ws2.Range("C8").Value = "4dd38fef-441a-4be0-b79a-ef6c9247ebcf.9ebd15eee1cc650407d011344150e433768ce247.zh-cn.xlf"
$ws2.Range("D8").Value = "2016-03-10 04:20:43"

Set-LinkText $ws2 '$A$6' "64fac99d-8603-4021-9534-a43a7f1cfd1d.md"
Set-LinkText $ws2 '$C$6' "64fac99d-8603-4021-9534-a43a7f1cfd1d.efd74baa64930741992ee27e4b2834cc8ef1667d.zh-cn.xlf"
Set-LinkText $ws2 '$A$7' "70058cc9-c335-4816-b6d1-eb41c2e484ea.md"
Set-LinkText $ws2 '$C$7' "70058cc9-c335-4816-b6d1-eb41c2e484ea.d5899fcb1515f857962642b1ceab8a68295ec2fb.zh-cn.xlf"
Set-LinkText $ws2 '$A$8' "4dd38fef-441a-4be0-b79a-ef6c9247ebcf.md"
Set-LinkText $ws2 '$C$8' "4dd38fef-441a-4be0-b79a-ef6c9247ebcf.9ebd15eee1cc650407d011344150e433768ce247.zh-cn.xlf"

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A6").Value = "64fac99d-8603-4021-9534-a43a7f1cfd1d.md"
$ws3.Range("B6").Value = "Handback transform failed"
$ws3.Range("C6").Value = "64fac99d-8603-4021-9534-a43a7f1cfd1d.efd74baa64930741992ee27e4b2834cc8ef1667d.de-de.xlf"
$ws3.Range("D6").Value = "2016-03-10 04:20:47"

$ws3.Range("A7").Value = "70058cc9-c335-4816-b6d1-eb41c2e484ea.md"
$ws3.Range("B7").Value = "In Translation"
$ws3.Range("C7").Value = "70058cc9-c335-4816-b6d1-eb41c2e484ea.d5899fcb1515f857962642b1ceab8a68295ec2fb.de-de.xlf"
$ws3.Range("D7").Value = "2016-03-10 04:13:23"

$ws3.Range("A8").Value = "4dd38fef-441a-4be0-b79a-ef6c9247ebcf.md"
$ws3.Range("B8").Value = "Ready for handoff"
$ws3.Range("C8").Value = "4dd38fef-441a-4be0-b79a-ef6c9247ebcf.9ebd15eee1cc650407d011344150e433768ce247.de-de.xlf"
$ws3.Range("D8").Value = "2016-03-10 04:20:47"

Set-LinkText $ws3 '$A$6' "64fac99d-8603-4021-9534-a43a7f1cfd1d.md"
Set-LinkText $ws3 '$C$6' "64fac99d-8603-4021-9534-a43a7f1cfd1d.efd74baa64930741992ee27e4b2834cc8ef1667d.de-de.xlf"
Set-LinkText $ws3 '$A$7' "70058cc9-c335-4816-b6d1-eb41c2e484ea.md"
Set-LinkText $ws3 '$C$7' "70058cc9-c335-4816-b6d1-eb41c2e484ea.d5899fcb1515f857962642b1ceab8a68295ec2fb.de-de.xlf"
Set-LinkText $ws3 '$A$8' "4dd38fef-441a-4be0-b79a-ef6c9247ebcf.md"
Set-LinkText $ws3 '$C$8' "4dd38fef-441a-4be0-b79a-ef6c9247ebcf.9ebd15eee1cc650407d011344150e433768ce247.de-de.xlf"
